$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (Panel A, 10yr Bond Futures) - Avg Daily Volume
$ws.Range("W11").Value = 38282.03402017841

# Row 26 (Panel B, E-mini Futures) - Avg Daily Volume
$ws.Range("D26").Value = 1372178.222222222
$ws.Range("E26").Value = 530511.9816199905
$ws.Range("G26").Value = 1233998
$ws.Range("H26").Value = 1647226.5
$ws.Range("I26").Value = 63
$ws.Range("J26").Value = 1438162.873015873
$ws.Range("K26").Value = 386608.0948163136
$ws.Range("L26").Value = 1156363
$ws.Range("M26").Value = 1385127
$ws.Range("N26").Value = 1650527
$ws.Range("O26").Value = 63
$ws.Range("P26").Value = 1591098.857142857
$ws.Range("Q26").Value = 384932.1306574436
$ws.Range("R26").Value = 1295451.5
$ws.Range("U26").Value = 63
$ws.Range("V26").Value = 1807131.333333333
$ws.Range("W26").Value = 483832.4142378493
$ws.Range("Z26").Value = 2037612
$ws.Range("AA26").Value = 63
$ws.Range("AB26").Value = 1722071.365079365
$ws.Range("AC26").Value = 494428.412564784
$ws.Range("AD26").Value = 1345607
$ws.Range("AE26").Value = 1656560
$ws.Range("AF26").Value = 2106191
$ws.Range("AG26").Value = 63

# Row 27 (Panel B, E-mini Futures) - Diff_Vol (Ann - Day)
$ws.Range("D27").Value = 218920.6349206349
$ws.Range("J27").Value = 152935.9841269841
$ws.Range("V27").Value = -216032.4761904762
$ws.Range("AB27").Value = -130972.5079365079

# Row 28 (Panel B, E-mini Futures) - # Obs
$ws.Range("D28").Value = 63
$ws.Range("J28").Value = 63
$ws.Range("P28").Value = 63
$ws.Range("V28").Value = 63
$ws.Range("AB28").Value = 63
